$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: the "_GoBack" bookmark that used to sit between
# "...benefit your organization" and ". Please contact..." is removed,
# and the two (identically formatted) runs around it collapse into a
# single run with the concatenated text.
# ------------------------------------------------------------------

$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

$merge = $d.Content.Duplicate
$merge.Find.Execute("organization. Please contact", $true, $false, $false, $false, $false, $true, 1, $false, "organization. Please contact", 2) | Out-Null

# ------------------------------------------------------------------
# Change 2: "North South University " -> "North South " and the
# "_GoBack" bookmark re-appears right after that run, collapsed
# (bookmarkStart immediately followed by bookmarkEnd, no content in
# between).
# ------------------------------------------------------------------

$uni = $d.Content.Duplicate
$uni.Find.Execute("University", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Re-create the "_GoBack" bookmark *before* deleting "University" (plus
# its trailing non-breaking space) so the anchor naturally settles at
# the end of the remaining "North South " text once the deletion
# happens (Word keeps bookmarks pinned to the surrounding text as it
# is edited).
$bmPoint = $d.Range($uni.Start, $uni.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

$uniFull = $d.Range($uni.Start, $uni.End + 1)
$uniFull.Delete()
